$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 data (6th measurement)
$ws.Range("C10").Value = -142358
$ws.Range("C10").Interior.Color = $ws.Range("C2").Interior.Color
$ws.Range("E10").Value = -362886
$ws.Range("E10").Interior.Color = $ws.Range("E2").Interior.Color
$ws.Range("G10").Formula = "=+E10-C10"
$ws.Range("H10").Value = 11000
$ws.Range("I10").Formula = "=G10/11"

# Move selection to I10 (also clears the stale topLeftCell scroll position)
[void]$ws.Range("I10").Select()
